$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price values that look numeric but are stored as text in the
# source workbook (e.g. "69.020.37" with multiple separators). Force text type
# via a temporary Text number format, then restore the original "Normal" style
# so no stray formatting is left behind on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.020.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.738.42"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.66%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.69"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.38%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +1.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.736.78"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.148"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.369"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.47%  "
$ws.Range("E12").Value = "  +1.66%  "
$ws.Range("E13").Value = "  -0.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.92"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.68%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.235.27"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.54%  "
$ws.Range("E16").Value = "  +2.88%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.009.61"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.720.61"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.07%  "
$ws.Range("E19").Value = "  +5.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "376.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E21").Value = "  +5.34%  "
$ws.Range("E22").Value = "  +3.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.02"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.26%  "
$ws.Range("E24").Value = "  +3.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.84%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.20"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.87%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000107"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.84%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "589.51"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.90%  "
$ws.Range("E31").Value = "  -0.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.43"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.82%  "
$ws.Range("E33").Value = "  +5.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.99"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.54%  "
$ws.Range("E35").Value = "  +4.89%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.63"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.94%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "162.67"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.16%  "
$ws.Range("E39").Value = "  +2.14%  "
$ws.Range("E40").Value = "  +4.10%  "
$ws.Range("E41").Value = "  +4.27%  "
$ws.Range("E42").Value = "  +3.83%  "
$ws.Range("E43").Value = "  +3.52%  "
$ws.Range("E44").Value = "  +1.22%  "
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "41.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.92%  "
$ws.Range("E47").Value = "  -2.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "156.26"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.12%  "
$ws.Range("E50").Value = "  +8.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.608"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.14%  "
